$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsOverview.Range("G2").Value = "2016-08-21 16:55:11"

$wsZhCn.Range("H2").Value = "2016-08-21 16:55:06"
$wsZhCn.Range("K2").Value = "2016-08-21 16:55:27"

$wsDeDe.Range("H2").Value = "2016-08-21 16:55:11"
$wsDeDe.Range("K2").Value = "2016-08-21 16:55:33"
